$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("AE2").Value = 3.048986434936523
$ws.Range("AF2").Value = 14.19855213165283
$ws.Range("AE3").Value = 3.158958435058594
$ws.Range("AF3").Value = 14.67528915405273
$ws.Range("AE4").Value = 3.253911733627319
$ws.Range("AF4").Value = 14.78826999664307
$ws.Range("AE5").Value = 3.339487314224243
$ws.Range("AF5").Value = 15.0087947845459
$ws.Range("AE6").Value = 3.268934011459351
$ws.Range("AF6").Value = 15.4818696975708
$ws.Range("AE7").Value = 2.643489360809326
$ws.Range("AF7").Value = 14.26635646820068
$ws.Range("AE8").Value = 2.766763687133789
$ws.Range("AF8").Value = 14.6733865737915
$ws.Range("AE9").Value = 2.947901725769043
$ws.Range("AF9").Value = 14.93699264526367
$ws.Range("AE10").Value = 3.03215479850769
$ws.Range("AF10").Value = 15.22404766082764
$ws.Range("AE11").Value = 3.107831478118896
$ws.Range("AF11").Value = 15.38049697875977
$ws.Range("AE12").Value = 3.298761129379272
$ws.Range("AF12").Value = 15.74327278137207
$ws.Range("AE13").Value = 3.354927778244019
$ws.Range("AF13").Value = 16.12627220153809
$ws.Range("AE14").Value = 1.762137293815613
$ws.Range("AF14").Value = 14.24727058410645
$ws.Range("AE15").Value = 2.164197683334351
$ws.Range("AF15").Value = 14.80567073822021
$ws.Range("AE16").Value = 2.394794702529907
$ws.Range("AF16").Value = 15.32405567169189
$ws.Range("AE17").Value = 2.70711612701416
$ws.Range("AF17").Value = 15.22656726837158
$ws.Range("AE18").Value = 2.857454061508179
$ws.Range("AF18").Value = 15.39282512664795
$ws.Range("AE19").Value = 3.000916719436646
$ws.Range("AF19").Value = 15.81495380401611
$ws.Range("AE20").Value = 3.196066379547119
$ws.Range("AF20").Value = 16.19318008422852
$ws.Range("AE21").Value = 3.363282442092896
$ws.Range("AF21").Value = 16.40904808044434
$ws.Range("AE22").Value = 3.431832075119019
$ws.Range("AF22").Value = 16.70996475219727
$ws.Range("AE23").Value = 1.478200554847717
$ws.Range("AF23").Value = 14.68713569641113
$ws.Range("AE24").Value = 1.718989968299866
$ws.Range("AF24").Value = 14.94417858123779
$ws.Range("AE25").Value = 2.011849880218506
$ws.Range("AF25").Value = 15.29533672332764
$ws.Range("AE26").Value = 2.246812582015991
$ws.Range("AF26").Value = 15.39813709259033
$ws.Range("AE27").Value = 2.483688354492188
$ws.Range("AF27").Value = 15.819580078125
$ws.Range("AE28").Value = 2.706009864807129
$ws.Range("AF28").Value = 16.15385818481445
$ws.Range("AE29").Value = 3.023434162139893
$ws.Range("AF29").Value = 16.43810844421387
$ws.Range("AE30").Value = 3.203835248947144
$ws.Range("AF30").Value = 16.62729835510254
$ws.Range("AE31").Value = 3.428520202636719
$ws.Range("AF31").Value = 16.89436721801758
$ws.Range("AE32").Value = 1.253238201141357
$ws.Range("AF32").Value = 14.86218929290771
$ws.Range("AE33").Value = 1.428144097328186
$ws.Range("AF33").Value = 15.21390724182129
$ws.Range("AE34").Value = 1.646275401115417
$ws.Range("AF34").Value = 15.43618869781494
$ws.Range("AE35").Value = 1.843245983123779
$ws.Range("AF35").Value = 15.77967262268066
$ws.Range("AE36").Value = 2.163870811462402
$ws.Range("AF36").Value = 16.05665588378906
$ws.Range("AE37").Value = 2.434584617614746
$ws.Range("AF37").Value = 16.35312843322754
$ws.Range("AE38").Value = 2.715826511383057
$ws.Range("AF38").Value = 16.59934234619141
$ws.Range("AE39").Value = 2.960421323776245
$ws.Range("AF39").Value = 16.86336326599121
$ws.Range("AE40").Value = 3.156507253646851
$ws.Range("AF40").Value = 17.02183532714844
$ws.Range("AE41").Value = 1.090214490890503
$ws.Range("AF41").Value = 15.0966968536377
$ws.Range("AE42").Value = 1.276663303375244
$ws.Range("AF42").Value = 15.32936000823975
$ws.Range("AE43").Value = 1.294190764427185
$ws.Range("AF43").Value = 15.5907564163208
$ws.Range("AE44").Value = 1.58684766292572
$ws.Range("AF44").Value = 15.97738075256348
$ws.Range("AE45").Value = 1.792916893959045
$ws.Range("AF45").Value = 16.12858200073242
$ws.Range("AE46").Value = 2.077682971954346
$ws.Range("AF46").Value = 16.52140235900879
$ws.Range("AE47").Value = 2.335130214691162
$ws.Range("AF47").Value = 16.83902740478516
$ws.Range("AE48").Value = 2.698145389556885
$ws.Range("AF48").Value = 16.9146785736084
$ws.Range("AE49").Value = 3.068830251693726
$ws.Range("AF49").Value = 17.04213523864746
$ws.Range("AE50").Value = 1.119083642959595
$ws.Range("AF50").Value = 15.52112483978271
$ws.Range("AE51").Value = 0.9675148129463196
$ws.Range("AF51").Value = 15.61314487457275
$ws.Range("AE52").Value = 1.146602749824524
$ws.Range("AF52").Value = 15.75729751586914
$ws.Range("AE53").Value = 1.300629377365112
$ws.Range("AF53").Value = 16.19294166564941
$ws.Range("AE54").Value = 1.578214526176453
$ws.Range("AF54").Value = 16.5205135345459
$ws.Range("AE55").Value = 1.819749593734741
$ws.Range("AF55").Value = 16.58085632324219
$ws.Range("AE56").Value = 2.028735399246216
$ws.Range("AF56").Value = 16.90520286560059
$ws.Range("AE57").Value = 2.44606876373291
$ws.Range("AF57").Value = 16.96755027770996
$ws.Range("AE58").Value = 2.717647552490234
$ws.Range("AF58").Value = 17.05212020874023
$ws.Range("AE59").Value = 1.166656136512756
$ws.Range("AF59").Value = 15.84199905395508
$ws.Range("AE60").Value = 0.9909321665763855
$ws.Range("AF60").Value = 16.08383369445801
$ws.Range("AE61").Value = 1.040380954742432
$ws.Range("AF61").Value = 16.37764549255371
$ws.Range("AE62").Value = 1.263261318206787
$ws.Range("AF62").Value = 16.6052360534668
$ws.Range("AE63").Value = 1.546685695648193
$ws.Range("AF63").Value = 16.79434394836426
$ws.Range("AE64").Value = 1.877916216850281
$ws.Range("AF64").Value = 16.81621551513672
$ws.Range("AE65").Value = 2.220521688461304
$ws.Range("AF65").Value = 16.89978408813477
$ws.Range("AE66").Value = 1.139411330223083
$ws.Range("AF66").Value = 16.07024192810059
$ws.Range("AE67").Value = 1.068029403686523
$ws.Range("AF67").Value = 16.1440486907959
$ws.Range("AE68").Value = 1.139615297317505
$ws.Range("AF68").Value = 16.52154922485352
$ws.Range("AE69").Value = 1.44741427898407
$ws.Range("AF69").Value = 16.44603729248047
$ws.Range("AE70").Value = 1.525559544563293
$ws.Range("AF70").Value = 16.81131744384766
